# Symbol-list refresh (GitHub Actions bot, 2022-12-27 14:38:53 UTC).
# "Price" (column D) values are stored as text in this sheet, so numeric-
# looking updates are written with a leading apostrophe to keep them as
# text, then the style is reset to "Normal" so no stray number-format /
# quote-prefix style sticks to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'24.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.383"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05901"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.399"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.501"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8121"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9375"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1422"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07396"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03097"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03085"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09342"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.871"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001573"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04723"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005993"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '17OneONEWorstin24h'

$ws.Range("D19").Value = "'0.005927"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001249"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.004749"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00008810"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.558"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'0.3222"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Value = "'0.0002654"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.03890"
$ws.Range("D40").Style = "Normal"

# Rows 41-43 rotate: KickToken -> row41, BKEXToken -> row42, CEJI -> row43.
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = "'0.006287"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '40KickTokenKICK'

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1070"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '41BKEXTokenBKK'

$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = "'0.002804"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '42CEJICEJI'

$ws.Range("D44").Value = "'0.008522"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005210"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.6714"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.001928"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("D50").Style = "Normal"
